# Model-Perf_Summary_v2.xlsx -- refresh the metrics table with the "Third Eval" /
# "Final Eval" models (finished the LDA/LS model round). The sheet is rebuilt from
# scratch (same as the upstream export script does) so the shared-string table is
# regenerated in the same column-major order: header row, then column B top-to-
# bottom, then column G top-to-bottom, which is what the new workbook stores.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Clear()

$nl = [char]10

# ---------------------------------------------------------------------------
# Header row (shared strings 0-5)
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Model"
$ws.Range("C1").Value = "Accuracy"
$ws.Range("D1").Value = "Precision"
$ws.Range("E1").Value = "Recall"
$ws.Range("F1").Value = "F1"
$ws.Range("G1").Value = "CM"

# ---------------------------------------------------------------------------
# Column B, top to bottom (new shared strings 6-17)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value  = "Random Forest-Feature Selection-Tune=Recall"
$ws.Range("B3").Value  = "Gradient Boost-FeatureSelection-Tune=Recall"
$ws.Range("B4").Value  = "Logistic Regression -- Second Eval - Untunned"
$ws.Range("B5").Value  = "Random Forest -- Second Eval - Tune=Recall"
$ws.Range("B6").Value  = "KNN -- Second Eval -- Tuning=Recall"
$ws.Range("B7").Value  = "LDA -- Second Eval - Tuning=Recall"
$ws.Range("B8").Value  = "Classification Tree -- Second Eval"
$ws.Range("B9").Value  = "GaussianNB"
$ws.Range("B10").Value = "Voting Classifier -- Second Eval --  RF(3)/GNB(1)/LDA(2)"
$ws.Range("B11").Value = "Logistic Regression -- Third Eval"
$ws.Range("B12").Value = "Random Forest -- Third Eval - w/oFeatures<0.1 -- Recall Tuning from Round 2"
$ws.Range("B13").Value = "Voting Classifier -- Final Eval --  RF3/LR3/GNB1"

# ---------------------------------------------------------------------------
# Column G, top to bottom (new shared strings 18-28)
# ---------------------------------------------------------------------------
$ws.Range("G2").Value  = "[[47 12]" + $nl + " [ 7  4]]"
$ws.Range("G3").Value  = "[[52  7]" + $nl + " [ 9  2]]"
$ws.Range("G4").Value  = "[[38 20]" + $nl + " [ 7  5]]"
$ws.Range("G5").Value  = "[[36 22]" + $nl + " [ 7  5]]"
$ws.Range("G6").Value  = "[[51  7]" + $nl + " [10  2]]"
$ws.Range("G7").Value  = "[[44 14]" + $nl + " [ 8  4]]"
$ws.Range("G8").Value  = "[[41 17]" + $nl + " [ 9  3]]"
$ws.Range("G9").Value  = "[[16 42]" + $nl + " [ 3  9]]"
$ws.Range("G10").Value = "[[30 28]" + $nl + " [ 5  7]]"
$ws.Range("G11").Value = "[[33 25]" + $nl + " [ 7  5]]"
$ws.Range("G12").Value = "[[33 25]" + $nl + " [ 7  5]]"
$ws.Range("G13").Value = "[[33 25]" + $nl + " [ 5  7]]"

# ---------------------------------------------------------------------------
# Column A (row index numbers, 0-based) -- plain numbers, no shared strings
# ---------------------------------------------------------------------------
For ($i = 0; $i -le 11; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

# ---------------------------------------------------------------------------
# Columns C-F, numeric metrics
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 0.7286
$ws.Range("D2").Value = 0.25
$ws.Range("E2").Value = 0.3636
$ws.Range("F2").Value = 0.2963

$ws.Range("C3").Value = 0.7714
$ws.Range("D3").Value = 0.2222
$ws.Range("E3").Value = 0.1818
$ws.Range("F3").Value = 0.2

$ws.Range("C4").Value = 0.6143
$ws.Range("D4").Value = 0.2
$ws.Range("E4").Value = 0.4167
$ws.Range("F4").Value = 0.2703

$ws.Range("C5").Value = 0.5857
$ws.Range("D5").Value = 0.1852
$ws.Range("E5").Value = 0.4167
$ws.Range("F5").Value = 0.2564

$ws.Range("C6").Value = 0.7571
$ws.Range("D6").Value = 0.2222
$ws.Range("E6").Value = 0.1667
$ws.Range("F6").Value = 0.1905

$ws.Range("C7").Value = 0.6857
$ws.Range("D7").Value = 0.2222
$ws.Range("E7").Value = 0.3333
$ws.Range("F7").Value = 0.2667

$ws.Range("C8").Value = 0.6286
$ws.Range("D8").Value = 0.15
$ws.Range("E8").Value = 0.25
$ws.Range("F8").Value = 0.1875

$ws.Range("C9").Value = 0.3571
$ws.Range("D9").Value = 0.1765
$ws.Range("E9").Value = 0.75
$ws.Range("F9").Value = 0.2857

$ws.Range("C10").Value = 0.5286
$ws.Range("D10").Value = 0.2
$ws.Range("E10").Value = 0.5833
$ws.Range("F10").Value = 0.2979

$ws.Range("C11").Value = 0.5429
$ws.Range("D11").Value = 0.1667
$ws.Range("E11").Value = 0.4167
$ws.Range("F11").Value = 0.2381

$ws.Range("C12").Value = 0.5429
$ws.Range("D12").Value = 0.1667
$ws.Range("E12").Value = 0.4167
$ws.Range("F12").Value = 0.2381

$ws.Range("C13").Value = 0.5714
$ws.Range("D13").Value = 0.2188
$ws.Range("E13").Value = 0.5833
$ws.Range("F13").Value = 0.3182

# ---------------------------------------------------------------------------
# Re-apply formatting: header row (B1:G1) and the index column (A2:A13) use the
# bold / thin-border / centered style from the original sheet.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$idxRange = $ws.Range("A2:A13")
$idxRange.Font.Bold = $true
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160
$idxRange.Borders.LineStyle = 1

# Column B is sized to fit the longest label.
$ws.Columns.Item(2).AutoFit()
